# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    9  = 2
    11 = -1
    24 = -4
    25 = -7
    26 = -2
    27 = 7
    29 = -8
    30 = 5
    34 = -6
    35 = -2
    37 = -7
    40 = 4
    43 = -2
    44 = -6
    48 = -4
    49 = -2
    51 = -4
    52 = -3
    54 = -1
    57 = -4
    66 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
